$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) — update "想去人数" (interested count) column F
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6896
$ws1.Range("F4").Value = 202
$ws1.Range("F5").Value = 44
$ws1.Range("F6").Value = 1070
$ws1.Range("F7").Value = 157

# Sheet "全部类型" (all types, aggregate of every category) — same rows, updated separately
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6897
$ws4.Range("F4").Value = 202
$ws4.Range("F5").Value = 44
$ws4.Range("F6").Value = 1070
$ws4.Range("F7").Value = 157
